# Adds a "LastImageEdit" tracking column to the Groups table and a
# "LastEdit" tracking column to the Expenses table (group/expense profile
# images now track when they were last uploaded/edited).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Groups table (rows 1-4) ---------------------------------------------
# Header I2: "LastEdit" -> "LastImageEdit"
$ws.Range("I2").Value = "LastImageEdit"

# Existing group "last edit" timestamps get new values (re-rolled ids)
$ws.Range("I3").Value = 551561551
$ws.Range("I4").Value = 645656454

# The blank filler rows below/around the Groups table no longer carry a
# (previously empty) cell in column I.
$ws.Range("I5").Clear()
$ws.Range("I6").Clear()
$ws.Range("I7").Clear()

# --- Expenses table (rows 8-10) -------------------------------------------
# A new "LastEdit" column is inserted at I, pushing the old "FK - GroupID"
# column (and its data) over to J.
$ws.Range("J8").Value = "FK - GroupID"          # header shifts right (was I8)
$ws.Range("I8").Value = "LastEdit"              # new header

$ws.Range("J9").Value = 1                       # old FK-GroupID value shifts right (was I9)
$ws.Range("I9").Value = 56465464                # new LastEdit value

$ws.Range("J10").Value = 2                      # old FK-GroupID value shifts right (was I10)
$ws.Range("I10").Value = 56151561               # new LastEdit value

# Give the new column a sensible display width (best-fit-ish).
$ws.Columns.Item(9).ColumnWidth = 11.6
